$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new product row: jose / 20000.0 / 40 / activo / Vinos
$ws.Range("A6").Value = "jose"

# "20000.0" must be stored as text (matches the Precio column convention
# used by the other rows), so force it via a leading apostrophe and then
# clear the resulting "quote prefix" style back to Normal so no extra
# cell style gets introduced.
$ws.Range("B6").Value = "'20000.0"
$ws.Range("B6").Style = "Normal"

$ws.Range("C6").Value = 40
$ws.Range("D6").Value = "activo"
$ws.Range("E6").Value = "Vinos"
